$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# --- Row 4 ---
$a4 = @'
Search the array
'@
$b4 = @'
lst=[ 1, 6, 3, 5, 3, 4 ] 
#checking if element 7 is present
# in the given list or not
i=7 
# if element present then return
# exist otherwise not exist
if i in lst: 
    print("exist") 
else: 
    print("not exist")
'@
$ws.Range("A4").Value = $a4
$ws.Range("B4").Value = $b4

# --- Row 5 ---
$a5 = @'
Max Consecutive Ones
'@
$b5 = @'
def search(input_list, num): 
    if (num in input_list):
        print("Element Found")
    else:
        print("Not Found")
search([12, 23, 45, 67, 6, 90] , 12)
'@
$ws.Range("A5").Value = $a5
$ws.Range("B5").Value = $b5

# --- Row 6 ---
$a6 = @'
Find Numbers with Even Number of Digits
'@
$b6 = @'
def sortedSquares(nums):
squares_list = []
for i in range(0, len(nums)):
square = nums[i] * nums[i];
squares_list.append(square)
sorted_squares_list = sorted(squares_list)
print sorted_squares_list;
return sorted_squares_list;
sortedSquares([-7,-3,2,3,11])   
'@
$ws.Range("A6").Value = $a6
$ws.Range("B6").Value = $b6

# --- Row 7 ---
$a7 = @'
Squares of a Sorted Array
'@
$b7 = @'
def sortedSquares(nums):
squares_list = []
for i in range(0, len(nums)):
square = nums[i] * nums[i];
squares_list.append(square)
\b
\b
sorted_squares_list = sorted(squares_list)
print sorted_squares_list;
return sorted_squares_list;
sortedSquares([-7,-3,2,3,11])
'@
$ws.Range("A7").Value = $a7
$ws.Range("B7").Value = $b7

# --- Formatting for B4:B7 (wrap text) ---
$ws.Range("B4:B7").WrapText = $true

# --- Special pasted-in font for B4:B5 (docs-Calibri, Google Docs paste artifact) ---
$ws.Range("B4:B5").Font.Name = "docs-Calibri"
$ws.Range("B4:B5").Font.Size = 11.25
$ws.Range("B4:B5").Font.Color = 0

# --- Row heights ---
$ws.Rows.Item(4).RowHeight = 142
$ws.Rows.Item(5).RowHeight = 102
$ws.Rows.Item(6).RowHeight = 151
$ws.Rows.Item(7).RowHeight = 129

# --- Column A width widened to fit new content ---
$ws.Columns.Item(1).ColumnWidth = 35.85

# --- Activate Sheet4 and select the last edited cell ---
$ws.Activate()
$ws.Range("B7").Select()
